$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.09 = 7903.59 pesos`n✅ 7903.59 pesos = 2.08 = 947.36 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the numeric rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("O10").Value = 3776.1
$ws2.Range("N12").Value = 3798
$ws2.Range("O12").Value = 455.245
